# Update "想去人数" (interest count) values in the F column across the
# four worksheets of the 杭州-漫展信息 workbook, per commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F7").Value  = 378
$ws.Range("F9").Value  = 69
$ws.Range("F11").Value = 668
$ws.Range("F12").Value = 1508
$ws.Range("F13").Value = 5867
$ws.Range("F15").Value = 1651
$ws.Range("F17").Value = 5572
$ws.Range("F22").Value = 1581
$ws.Range("F23").Value = 828
$ws.Range("F26").Value = 1174
$ws.Range("F28").Value = 159
$ws.Range("F29").Value = 15
$ws.Range("F31").Value = 3834

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 105
$ws.Range("F5").Value = 203
$ws.Range("F8").Value = 315

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9456
$ws.Range("F3").Value = 589
$ws.Range("F4").Value = 2183
$ws.Range("F5").Value = 546

# Sheet "全部类型" (All types - combined view)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 9456
$ws.Range("F3").Value  = 589
$ws.Range("F4").Value  = 2183
$ws.Range("F7").Value  = 546
$ws.Range("F9").Value  = 378
$ws.Range("F11").Value = 69
$ws.Range("F14").Value = 668
$ws.Range("F15").Value = 1508
$ws.Range("F16").Value = 5867
$ws.Range("F18").Value = 315
$ws.Range("F19").Value = 1651
$ws.Range("F25").Value = 5572
$ws.Range("F30").Value = 1581
$ws.Range("F31").Value = 828
$ws.Range("F34").Value = 1174
$ws.Range("F36").Value = 159
$ws.Range("F40").Value = 15
$ws.Range("F45").Value = 3834
